# Fill in the 5th person's (Yasaman Tavakoli, row 10) quiz scores that were
# previously left blank, clear the "still needs data" highlight color from
# that row, and move the active selection to the next entry cell (K11),
# matching the author's "5 person report added yasaman last person in
# tuesday" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q1 HW1, Q1 Seminar, Q2 Circle, Q3 Factory, Q4 Complex, Q5 Library scores
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 88
$ws.Range("J10").Value = 70
$ws.Range("K10").Value = 100

# The row was highlighted (orange accent fill) to mark it as incomplete;
# now that it is filled in, clear that highlight back to white/background.
$ws.Range("A10:K10").Interior.ThemeColor = 2
$ws.Range("A10:K10").Interior.TintAndShade = 0

# Move the active cell/selection to the next row's entry point.
$ws.Range("K11").Select() | Out-Null
